# "Corrected code of Parser" - adds a few more parser trace rows/cells to
# Sheet1 (rows 35-38), plus a handful of fixes to existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- Row 9: new cell E9 ---
$ws.Range("E9").Value = "S34"

# --- Row 18: new cells B18, U18 ---
$ws.Range("B18").Value = "S23"
$ws.Range("U18").Value = 22

# --- Row 19: U19 value corrected ---
$ws.Range("U19").Value = 25

# --- Row 27: N27 value corrected ---
$ws.Range("N27").Value = "R11"

# --- Row 33: new cell I33 ---
$ws.Range("I33").Value = "R6"

# --- Row 35 (new row) ---
$ws.Range("A35").Value = 33
$ws.Range("I35").Value = "R2"
$ws.Range("M35").Value = "R2"

# --- Row 36 (new row) ---
$ws.Range("A36").Value = 34
$ws.Range("B36").Value = "S12"
$ws.Range("E36").Value = "S34"
$ws.Range("S36").Value = 35
$ws.Range("T36").Value = 10
$ws.Range("U36").Value = 11

# --- Row 37 (new row) ---
$ws.Range("A37").Value = 35
$ws.Range("F37").Value = "S36"
$ws.Range("H37").Value = "S16"
$ws.Range("J37").Value = "S15"

# --- Row 38 (new row) ---
$ws.Range("A38").Value = 36
$ws.Range("F38").Value = "R14"
$ws.Range("H38").Value = "R14"
$ws.Range("J38").Value = "R14"
$ws.Range("K38").Value = "R14"
$ws.Range("L38").Value = "R14"
$ws.Range("N38").Value = "R14"

# --- View state: scroll down to the new rows and select F37 ---
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F37").Select()
